$d = $word.ActiveDocument

# Remove the "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph and the
# "© 2020 . Contact: luizeleno@usp.br. ..." footer paragraph (and the blank
# paragraph that sat between "LOQ4038: ..." and the footer block) that the
# Jekyll site build no longer emits.

$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*Ver no Jupiter Salvar em pdf Salvar em docx*" -or $t -like "*Contact: luizeleno@usp.br*") {
        $p.Range.Delete()
    }
}

$count2 = $d.Paragraphs.Count
for ($i = $count2; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*LOQ4038*") {
        $next = $p.Next()
        $next.Range.Delete()
        break
    }
}
